# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in AC1:AE1
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the format from an existing header cell.
$ws.Range("Z1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record for every player row (2-37).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 58
    $ws.Cells.Item($r, 30).Value = 56
    $ws.Cells.Item($r, 31).Value = 0
}
